# Add a second test case (CreateAccountTest) to the test data workbook,
# mirroring the existing LoginTest / test_suite pattern used for POM
# (Page Object Model) based Selenium tests.
#
# NOTE: worksheet references are re-fetched (via Worksheets.Item(name))
# after any operation that reorders sheets (Add / Move), since stale
# references can end up pointing at the wrong physical sheet afterwards.

$wb = $excel.ActiveWorkbook

# 1. Create the new worksheet (left in place for now; it gets moved to
#    the end once all sheets/cells have been populated).
$newSheet = $wb.Worksheets.Add()
$newSheet.Name = "CreateAccountTest"

# 2. Register the new test case as an additional row in the test_suite
#    sheet (TCID / Runmode columns).
$testSuite = $wb.Worksheets.Item("test_suite")
$testSuite.Range("A3").Value = "CreateAccountTest"
$testSuite.Range("B3").Value = "Y"

# 3. Populate the new CreateAccountTest sheet with its test data.
$newSheet.Range("A1").Value = "accountname"
$newSheet.Range("A2").Value = "george"

# 4. Move the new sheet so it becomes the last sheet (after LoginTest).
$loginTest = $wb.Worksheets.Item("LoginTest")
$createAccountTest = $wb.Worksheets.Item("CreateAccountTest")
$createAccountTest.Move($null, $loginTest)

# 5. Restore/update the selections on each sheet, re-fetching sheet
#    references fresh (post-move) so the right physical sheet is touched,
#    and make CreateAccountTest the active tab.
$testSuite = $wb.Worksheets.Item("test_suite")
$testSuite.Range("B3").Select()

$createAccountTest = $wb.Worksheets.Item("CreateAccountTest")
$createAccountTest.Activate()
$createAccountTest.Range("A2").Select()

$wb.Save()
